$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "56.919.44"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "3.238.58"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "396.26"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.56"
$ws.Range("E6").Value = "  -3.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.581"
$ws.Range("E7").Value = "  +4.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -1.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.14"
$ws.Range("E10").Value = "  -1.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0949"
$ws.Range("E11").Value = "  +5.02%  "
$ws.Range("E12").Value = "  +1.93%  "
$ws.Range("D13").Value = "3.748.34"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.27"
$ws.Range("E14").Value = "  +2.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "18.83"
$ws.Range("E15").Value = "  -2.07%  "
$ws.Range("D16").Value = "3.244.98"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("E17").Value = "  -4.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.95"
$ws.Range("E18").Value = "  +2.26%  "
$ws.Range("D19").Value = "56.690.37"
$ws.Range("E19").Value = "  +0.74%  "
$ws.Range("E20").Value = "  -2.94%  "
$ws.Range("E21").Value = "  +6.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.82"
$ws.Range("E22").Value = "  -2.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "291.41"
$ws.Range("E23").Value = "  -0.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.15"
$ws.Range("E24").Value = "  -0.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.17"
$ws.Range("E25").Value = "  -2.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.96"
$ws.Range("E26").Value = "  -2.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.99"
$ws.Range("E27").Value = "  -0.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.37"
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.169"
$ws.Range("E29").Value = "  -1.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.22"
$ws.Range("E30").Value = "  -4.04%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  -3.00%  "
$ws.Range("E33").Value = "  -1.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.96"
$ws.Range("E34").Value = "  +10.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0482"
$ws.Range("E35").Value = "  -2.84%  "
$ws.Range("E36").Value = "  +0.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.11"
$ws.Range("E37").Value = "  -0.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("E39").Value = "  -3.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.98"
$ws.Range("E40").Value = "  -4.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "136.80"
$ws.Range("E41").Value = "  +0.21%  "
$ws.Range("E42").Value = "  +1.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.91"
$ws.Range("E43").Value = "  -3.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.86"
$ws.Range("E44").Value = "  -3.65%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.282"
$ws.Range("E45").Value = "  -1.45%  "
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.59"
$ws.Range("E46").Value = "  -3.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.23"
$ws.Range("E47").Value = "  +6.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.28"
$ws.Range("E48").Value = "  -1.82%  "
$ws.Range("D49").Value = "2.150.39"
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("E50").Value = "  -4.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.95"
$ws.Range("E51").Value = "  -9.21%  "
